# Update version string across the workbook for the new release:
#   "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
# becomes
#   "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet updates ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: " + $newVersion

$about.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Pansan Coal Mine, China, M5206, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet updates ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

# S2:S13 hold the per-row build_version text
for ($row = 2; $row -le 13; $row++) {
    $cell = $data.Cells.Item($row, 19)  # column S = 19
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
